# Apply dual-variable updates for "Working version of branch and price with L-shaped subproblems"
$wb = $excel.ActiveWorkbook

# --- Sheet: u_MAB ---
$ws = $wb.Worksheets.Item("u_MAB")
$ws.Range("B15").Value = 0
$ws.Range("A16").Value = 0
$ws.Range("A24").Value = 1.364433831044545
$ws.Range("A27").Value = 0.3976231535489507
$ws.Range("A40").Value = 0
$ws.Range("B40").Value = 0
$ws.Range("B48").Value = 0.94746437695053
$ws.Range("B49").Value = 0.5079237467546113
$ws.Range("A51").Value = 0.2061950664522399
$ws.Range("A52").Value = 0.05182702263477305
$ws.Range("B61").Value = 0

# --- Sheet: u_EOH ---
$ws = $wb.Worksheets.Item("u_EOH")
$ws.Range("A2").Value = -0.3203540442204794
$ws.Range("A3").Value = -0.2191222914782864

# --- Sheet: v_l ---
$ws = $wb.Worksheets.Item("v_l")
$ws.Range("A2").Value = 2229019.522400185
$ws.Range("A3").Value = 1994659.418273741
$ws.Range("A4").Value = 9323831.392055079

$wb.Save()
